$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 450
$ws.Range("C6").Select()
